$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6757276
$ws.Range("F2").Value = "BSC Rapid Chemnitz"
$ws.Range("G2").Value = "FV Dresden 06 Laubegast"
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 2.6
$ws.Range("M2").Value = 2.25
$ws.Range("N2").Value = 2.6
$ws.Range("O2").Value = 3.5
$ws.Range("P2").Value = 2.25
$ws.Range("Q2").Value = 0.25
$ws.Range("R2").Value = 1.75
$ws.Range("S2").Value = 2.05
$ws.Range("T2").Value = 3.25
$ws.Range("U2").Value = 1.775
$ws.Range("V2").Value = 2.025
$ws.Range("W2").Value = 1.6
$ws.Range("Z2").Value = 0.75
$ws.Range("AB2").Value = -0.5
$ws.Range("AC2").Value = 0.5125
$ws.Range("B3").Value = 6760228
$ws.Range("F3").Value = "SpVg Porz 1919"
$ws.Range("G3").Value = "Bonn Endenich 1908"
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 1
$ws.Range("K3").Value = 2.25
$ws.Range("M3").Value = 2.625
$ws.Range("N3").Value = 1.75
$ws.Range("O3").Value = 3.6
$ws.Range("P3").Value = 3.8
$ws.Range("Q3").Value = -0.5
$ws.Range("R3").Value = 1.8
$ws.Range("S3").Value = 2
$ws.Range("T3").Value = 2.75
$ws.Range("U3").Value = 1.8
$ws.Range("V3").Value = 2
$ws.Range("W3").Value = 0.75
$ws.Range("Z3").Value = 0.8
$ws.Range("AB3").Value = 0.8
$ws.Range("AC3").Value = -1
$ws.Range("B4").Value = 6781316
$ws.Range("F4").Value = "SV Schott Jena"
$ws.Range("G4").Value = "SV 09 Arnstadt"
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = "A"
$ws.Range("K4").Value = 5
$ws.Range("L4").Value = 4.5
$ws.Range("M4").Value = 1.45
$ws.Range("N4").Value = 6.5
$ws.Range("O4").Value = 4.333
$ws.Range("P4").Value = 1.363
$ws.Range("Q4").Value = 1.5
$ws.Range("R4").Value = 1.825
$ws.Range("S4").Value = 1.975
$ws.Range("T4").Value = 3
$ws.Range("U4").Value = 1.825
$ws.Range("V4").Value = 1.975
$ws.Range("W4").Value = -1
$ws.Range("Y4").Value = 0.363
$ws.Range("Z4").Value = -1
$ws.Range("AA4").Value = 0.9750000000000001
$ws.Range("AC4").Value = 0.9750000000000001
$ws.Range("B5").Value = 6781315
$ws.Range("F5").Value = "SSV Markranstadt"
$ws.Range("G5").Value = "BSC Rapid Chemnitz"
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = "H"
$ws.Range("K5").Value = 1.25
$ws.Range("L5").Value = 4.75
$ws.Range("M5").Value = 10
$ws.Range("N5").Value = 1.222
$ws.Range("O5").Value = 5.25
$ws.Range("P5").Value = 8.5
$ws.Range("Q5").Value = -2
$ws.Range("R5").Value = 1.925
$ws.Range("S5").Value = 1.875
$ws.Range("T5").Value = 3.5
$ws.Range("U5").Value = 1.775
$ws.Range("V5").Value = 1.925
$ws.Range("W5").Value = 0.222
$ws.Range("Y5").Value = -1
$ws.Range("Z5").Value = 0
$ws.Range("AA5").Value = -0
$ws.Range("AC5").Value = 0.925
$ws.Range("B11").Value = 7035048
$ws.Range("F11").Value = "SG Unterrath"
$ws.Range("G11").Value = "TuRU Dsseldorf"
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "H"
$ws.Range("K11").Value = 3.25
$ws.Range("L11").Value = 4
$ws.Range("M11").Value = 1.8
$ws.Range("N11").Value = 2.9
$ws.Range("O11").Value = 4
$ws.Range("P11").Value = 1.95
$ws.Range("Q11").Value = 0.5
$ws.Range("T11").Value = 3
$ws.Range("U11").Value = 1.75
$ws.Range("V11").Value = 1.95
$ws.Range("W11").Value = 1.9
$ws.Range("Y11").Value = -1
$ws.Range("Z11").Value = 0.8
$ws.Range("AA11").Value = -1
$ws.Range("AC11").Value = 0.95
$ws.Range("B12").Value = 7035046
$ws.Range("F12").Value = "Cronenberger SC"
$ws.Range("G12").Value = "FC Viersen"
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 2
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 3.6
$ws.Range("M12").Value = 3
$ws.Range("N12").Value = 2
$ws.Range("P12").Value = 3
$ws.Range("R12").Value = 1.8
$ws.Range("S12").Value = 2
$ws.Range("T12").Value = 2.75
$ws.Range("U12").Value = 1.8
$ws.Range("V12").Value = 2
$ws.Range("Y12").Value = 2
$ws.Range("AA12").Value = 1
$ws.Range("AB12").Value = -1
$ws.Range("AC12").Value = 1
$ws.Range("B13").Value = 7035047
$ws.Range("F13").Value = "SC Dsseldorf West"
$ws.Range("G13").Value = "VfL Viktoria JuchenGarzweiler"
$ws.Range("H13").Value = 3
$ws.Range("I13").Value = 4
$ws.Range("J13").Value = "A"
$ws.Range("K13").Value = 1.909
$ws.Range("L13").Value = 3.75
$ws.Range("M13").Value = 3.1
$ws.Range("N13").Value = 2.2
$ws.Range("O13").Value = 3.6
$ws.Range("P13").Value = 2.625
$ws.Range("Q13").Value = -0.25
$ws.Range("R13").Value = 2
$ws.Range("S13").Value = 1.8
$ws.Range("U13").Value = 1.825
$ws.Range("V13").Value = 1.975
$ws.Range("W13").Value = -1
$ws.Range("Y13").Value = 1.625
$ws.Range("Z13").Value = -1
$ws.Range("AA13").Value = 0.8
$ws.Range("AB13").Value = 0.825
$ws.Range("AC13").Value = -1
$ws.Range("F21").Value = "SG Unterrath"
$ws.Range("G25").Value = "SG Unterrath"
$ws.Range("F39").Value = "SSV Markranstadt"
$ws.Range("F55").Value = "SC Dsseldorf West"
$ws.Range("F61").Value = "SpVg Porz 1919"
$ws.Range("F62").Value = "Cronenberger SC"
$ws.Range("F65").Value = "SV Schott Jena"
$ws.Range("F69").Value = "SG Unterrath"
$ws.Range("G71").Value = "SG Unterrath"

Write-Output "Applied changes"
